# Add ui and ux
# New columns: F = id, G = createdAt
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "id"
$ws.Range("G1").Value = "createdAt"

# New member rows appended below the existing data (rows 3-6).
# Columns B (phone) and F (id) hold digit-only strings that Excel would
# otherwise auto-coerce to numbers, so force text format before writing.
$newRows = @(
    @{ Row=3; Name="Abu Inshah"; Phone="7449085120"; Email="ajai17101999@gmail.com";      Designation="Health insurance advisor"; Photo="uploads/abu_inshah_1752306015359.jpeg"; Id=$null;          CreatedAt=$null },
    @{ Row=4; Name="ram";        Phone="7449085120"; Email="aiautomationhig@gmail.com";   Designation="Health insurance advisor"; Photo="uploads/ram_1752311161933.jpeg";        Id="1752311162001"; CreatedAt="2025-07-12T09:06:02.001Z" },
    @{ Row=5; Name="ram";        Phone="7449085120"; Email="selvasuresh460@gmail.com";    Designation="Health insurance advisor"; Photo="uploads/ram_1752314093239.jpeg";        Id="1752314093256"; CreatedAt="2025-07-12T09:54:53.256Z" },
    @{ Row=6; Name="Abu Inshah"; Phone="7449085120"; Email="wealthplusacademy@gmail.com"; Designation="Wealth Manager";            Photo="uploads/abu_inshah_1752314719383.jpeg"; Id="1752314719399"; CreatedAt="2025-07-12T10:05:19.399Z" }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Range("A$row").Value = $r.Name

    $ws.Range("B$row").NumberFormat = "@"
    $ws.Range("B$row").Value = $r.Phone
    $ws.Range("B$row").Style = "Normal"

    $ws.Range("C$row").Value = $r.Email
    $ws.Range("D$row").Value = $r.Designation
    $ws.Range("E$row").Value = $r.Photo

    if ($r.Id -ne $null) {
        $ws.Range("F$row").NumberFormat = "@"
        $ws.Range("F$row").Value = $r.Id
        $ws.Range("F$row").Style = "Normal"
    }

    if ($r.CreatedAt -ne $null) {
        $ws.Range("G$row").Value = $r.CreatedAt
    }
}
